# Drop.xlsx - "add some new veges"
# The four separate resource-gathering drop rows (水池/Pool, 矿洞/Mine,
# 蘑菇/Mushroom, 枯木/Deadwood) are removed, and the surviving row (originally
# "草丛"/Grass, id 23000001) is turned into a new "丝瓜"/Loofah drop entry
# listing the new vegetable materials. Every row below shifts up accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 4 rows that get folded away (水池, 矿洞, 蘑菇, 枯木).
# This auto-shifts every following row up by 4 and resizes the table/filter.
$ws.Rows("5:8").Delete()

# Turn the remaining first data row into the new "丝瓜" (loofah) entry.
$ws.Range("C4").Value = "dlsigua"
$ws.Range("B4").Value = "丝瓜"
$ws.Range("D4").Value = "zzwandou;zzyumi;zzpingguo;zzlanmei;zznangua;zzxihongshi;zzqiezi;zzluobo;zztudou;zzlajiao;zzyangcong"
$ws.Range("E4").Value = "9;9;9;9;9;9;9;9;9;9;9"

# Those cells were highlighted red (marked for edit); clear that highlight now
# that the row has its final content.
$ws.Range("B4:D4").Interior.Pattern = -4142

# Leave the selection where the author last left it.
$ws.Range("E4").Select()
